$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row per the refreshed snapshot
$ws.Range("D2").Value = "'27.719.43"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "'1.878.35"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'332.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4710"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("D8").Value = "'0.3952"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("D9").Value = "'47.84"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'0.08048"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.65%  "
$ws.Range("D11").Value = "'1.029"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").Value = "'22.15"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("D13").Value = "'1.880.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.57%  "
$ws.Range("D14").Value = "'5.980"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.44%  "
$ws.Range("D15").Value = "'7.126"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'0.00001048"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "'87.10"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "'0.06674"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("D20").Value = "'17.24"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'27.723.79"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "'5.532"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("E24").Value = "  +1.36%  "
$ws.Range("D25").Value = "'2.309"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "'2.099.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").Value = "'159.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'20.19"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").Value = "'2.103"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("E30").Value = "  +2.71%  "
$ws.Range("D31").Value = "'121.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").Value = "'0.9819"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("D33").Value = "'0.09513"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("D34").Value = "'1.445"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("D35").Value = "'3.595"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "'5.368"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").Value = "'0.06123"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.81%  "
$ws.Range("D38").Value = "'0.02261"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "'1.230"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("D40").Value = "'8.099"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "'0.6018"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'0.1900"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.21%  "
$ws.Range("D43").Value = "'10.35"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "'1.267"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").Value = "'0.5726"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.52%  "
$ws.Range("D46").Value = "'12.24"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.42%  "
$ws.Range("E47").Value = "  +2.23%  "
$ws.Range("D48").Value = "'3.389"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.44%  "
$ws.Range("D49").Value = "'0.06912"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "'114.45"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.83%  "

# Row 51: EOS dropped out of the top list, replaced by BabyDogeCoin
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000300"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.21%  "
